$d = $word.ActiveDocument

# Update the date/day heading (unique text in the document).
$d.Content.Find.Execute("2024-03-12 Tuesday", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2024-03-13 Wednesday", 2)

# Update the division problems in the grid. Several cells share identical
# source text (e.g. "18÷3=" and "59÷3=" each appear twice), so address each
# cell directly via the Tables object model rather than a global find/replace.
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "78÷9="
$t.Cell(1, 2).Range.Text = "97÷2="
$t.Cell(1, 3).Range.Text = "26÷3="
$t.Cell(1, 4).Range.Text = "85÷8="
$t.Cell(1, 5).Range.Text = "11÷7="

$t.Cell(5, 1).Range.Text = "76÷5="
$t.Cell(5, 2).Range.Text = "10÷8="
$t.Cell(5, 3).Range.Text = "89÷6="
$t.Cell(5, 4).Range.Text = "26÷6="
$t.Cell(5, 5).Range.Text = "75÷5="

$t.Cell(9, 1).Range.Text = "15÷6="
$t.Cell(9, 2).Range.Text = "27÷2="
$t.Cell(9, 3).Range.Text = "43÷3="
$t.Cell(9, 4).Range.Text = "25÷3="
$t.Cell(9, 5).Range.Text = "55÷4="

$t.Cell(13, 1).Range.Text = "82÷9="
$t.Cell(13, 2).Range.Text = "99÷9="
$t.Cell(13, 3).Range.Text = "14÷4="
$t.Cell(13, 4).Range.Text = "82÷4="
$t.Cell(13, 5).Range.Text = "74÷7="

$t.Cell(17, 1).Range.Text = "45÷9="
$t.Cell(17, 2).Range.Text = "93÷4="
$t.Cell(17, 3).Range.Text = "73÷8="
$t.Cell(17, 4).Range.Text = "51÷2="
$t.Cell(17, 5).Range.Text = "80÷7="
